$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.702.41'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '3.439.65'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.23'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.04'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.442.12'
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  +8.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.34'
$ws.Range("E10").Value = '  -2.99%  '
$ws.Range("E11").Value = '  +2.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.442'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").Value = '4.037.98'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("E14").Value = '  -2.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000194'
$ws.Range("E15").Value = '  +3.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.29'
$ws.Range("E16").Value = '  +3.65%  '
$ws.Range("D17").Value = '64.731.40'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '3.480.21'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.36'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.28'
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '386.10'
$ws.Range("E21").Value = '  -1.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.17'
$ws.Range("E22").Value = '  -3.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.27'
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.545'
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000124'
$ws.Range("E26").Value = '  +14.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.80'
$ws.Range("E27").Value = '  +2.74%  '
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.18'
$ws.Range("E30").Value = '  +6.65%  '
$ws.Range("E31").Value = '  +4.07%  '
$ws.Range("E32").Value = '  -0.10%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.55'
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.63'
$ws.Range("E34").Value = '  +0.39%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.08'
$ws.Range("E36").Value = '  +3.35%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.51'
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '162.95'
$ws.Range("E38").Value = '  +2.88%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.017.41'
$ws.Range("E39").Value = '  +5.17%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("E40").Value = '  +1.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0765'
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.24'
$ws.Range("E42").Value = '  -3.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.56'
$ws.Range("E43").Value = '  +4.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.84'
$ws.Range("E44").Value = '  +2.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0316'
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.773'
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.70'
$ws.Range("E47").Value = '  +9.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.08'
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.880'
$ws.Range("E49").Value = '  +6.59%  '
$ws.Range("E50").Value = '  +3.55%  '
$ws.Range("E51").Value = '  +4.21%  '
